$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. settings sheet: rename the "form_id" setting to "table_id"
#    (same row/position, just a relabelled key -- the sharedStrings
#    table will naturally drop "form_id" and gain "table_id")
# ------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Cells.Item(2, 1).Value = "table_id"

# ------------------------------------------------------------------
# 2. add a new "properties" worksheet at the end of the workbook
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$props = $wb.Worksheets.Add($null, $lastSheet)
$props.Name = "properties"

# colOrder value (json list of column names) -- entered first
$props.Cells.Item(2, 5).Value = '["FOL_date","FOL_B_AnimID","FOL_CL_community_id","FOL_time_begin","FOL_time_end","FOL_flag_begin_in_nest","FOL_flag_end_in_nest","FOL_duration","FOL_distance_traveled","FOL_am_observer1","FOL_am_observer2","FOL_pm_observer1","FOL_pm_observer2","FOL_study_code1","FOL_study_code2","FOL_day","FOL_month","FOL_year"]'

# header row
$props.Cells.Item(1, 1).Value = "partition"
$props.Cells.Item(1, 2).Value = "aspect"
$props.Cells.Item(1, 3).Value = "key"
$props.Cells.Item(1, 4).Value = "type"
$props.Cells.Item(1, 5).Value = "value"

# colOrder
$props.Cells.Item(2, 1).Value = "Table"
$props.Cells.Item(2, 2).Value = "default"
$props.Cells.Item(2, 3).Value = "colOrder"
$props.Cells.Item(2, 4).Value = "array"

# defaultViewType
$props.Cells.Item(3, 1).Value = "Table"
$props.Cells.Item(3, 2).Value = "default"
$props.Cells.Item(3, 3).Value = "defaultViewType"
$props.Cells.Item(3, 4).Value = "string"
$props.Cells.Item(3, 5).Value = "LIST"

# detailViewFileName
$props.Cells.Item(4, 1).Value = "Table"
$props.Cells.Item(4, 2).Value = "default"
$props.Cells.Item(4, 3).Value = "detailViewFileName"
$props.Cells.Item(4, 4).Value = "configpath"
$props.Cells.Item(4, 5).Value = "config/tables/follow/html/follow_detail.html"

# listViewFileName
$props.Cells.Item(5, 1).Value = "Table"
$props.Cells.Item(5, 2).Value = "default"
$props.Cells.Item(5, 3).Value = "listViewFileName"
$props.Cells.Item(5, 4).Value = "configpath"
$props.Cells.Item(5, 5).Value = "config/tables/follow/html/follow_list.html"

# ------------------------------------------------------------------
# 3. selections: settings -> A3, properties (new, now active) -> C10
# ------------------------------------------------------------------
$settings.Activate()
$settings.Range("A3").Select()

$props.Activate()
$props.Range("C10").Select()
